# Append new Streamlit dashboard code block after the existing content,
# matching the authored diff: a horizontal rule (paragraph with a bottom
# border), then a series of "# N. ..." comment / st.subheader / st.plotly_chart
# paragraphs (separated by blank paragraphs), ending with the 10B block.

$d = $word.ActiveDocument
$sel = $word.Selection

$items = @(
    @{ Kind = "BORDER"; Text = "" },
    @{ Kind = "EMPTY"; Text = "" },
    @{ Kind = "TEXT"; Text = "# 1. Total vs Average Call Days by Division" },
    @{ Kind = "TEXT"; Text = "st.subheader(`"1. Total vs Average Call Days by Division`")" },
    @{ Kind = "TEXT"; Text = "st.plotly_chart(fig_call_days, use_container_width=True)" },
    @{ Kind = "EMPTY"; Text = "" },
    @{ Kind = "TEXT"; Text = "# 2. Doctor Call Average by Division" },
    @{ Kind = "TEXT"; Text = "st.subheader(`"2. Doctor Call Avg by Division`")" },
    @{ Kind = "TEXT"; Text = "st.plotly_chart(fig_doc_avg, use_container_width=True)" },
    @{ Kind = "EMPTY"; Text = "" },
    @{ Kind = "TEXT"; Text = "# 3. Plan vs Actual DR Calls" },
    @{ Kind = "TEXT"; Text = "st.subheader(`"3. Plan vs Actual DR Calls`")" },
    @{ Kind = "TEXT"; Text = "st.plotly_chart(fig_plan_actual, use_container_width=True)" },
    @{ Kind = "EMPTY"; Text = "" },
    @{ Kind = "TEXT"; Text = "# 4. 2PC Frequency Coverage %" },
    @{ Kind = "TEXT"; Text = "st.subheader(`"4. 2PC Frequency Coverage % by Division`")" },
    @{ Kind = "TEXT"; Text = "st.plotly_chart(fig_2pc, use_container_width=True)" },
    @{ Kind = "EMPTY"; Text = "" },
    @{ Kind = "TEXT"; Text = "# 5. Total DR Coverage %" },
    @{ Kind = "TEXT"; Text = "st.subheader(`"5. Total DR Coverage % by Division`")" },
    @{ Kind = "TEXT"; Text = "st.plotly_chart(fig_total_cov, use_container_width=True)" },
    @{ Kind = "EMPTY"; Text = "" },
    @{ Kind = "TEXT_PAGEBREAK"; Text = "# 6. Field Work, Leaves, Total Days Comparison" },
    @{ Kind = "TEXT"; Text = "st.subheader(`"6. Comparison of Working Days by Division`")" },
    @{ Kind = "TEXT"; Text = "st.plotly_chart(fig_working_days, use_container_width=True)" },
    @{ Kind = "EMPTY"; Text = "" },
    @{ Kind = "TEXT"; Text = "# 7. Call and Visit Trends by Zone" },
    @{ Kind = "TEXT"; Text = "st.subheader(`"7. Calls Trends by Zone`")" },
    @{ Kind = "TEXT"; Text = "st.plotly_chart(fig_zone_trend, use_container_width=True)" },
    @{ Kind = "EMPTY"; Text = "" },
    @{ Kind = "TEXT"; Text = "# 8. Gauge Charts for 1PC and 2PC Coverage %" },
    @{ Kind = "TEXT"; Text = "st.subheader(`"8. Coverage Gauges (1PC & 2PC)`")" },
    @{ Kind = "TEXT"; Text = "col7, col8 = st.columns(2)" },
    @{ Kind = "TEXT"; Text = "with col7:" },
    @{ Kind = "TEXT"; Text = "    st.plotly_chart(fig1pc)" },
    @{ Kind = "TEXT"; Text = "with col8:" },
    @{ Kind = "TEXT"; Text = "    st.plotly_chart(fig2pc)" },
    @{ Kind = "EMPTY"; Text = "" },
    @{ Kind = "TEXT"; Text = "# 9. Doctor Visit Distribution (Total / Visited / Missed)" },
    @{ Kind = "TEXT"; Text = "st.subheader(`"9. Doctor Visit Distribution`")" },
    @{ Kind = "TEXT"; Text = "st.plotly_chart(fig_visit_dist, use_container_width=True)" },
    @{ Kind = "EMPTY"; Text = "" },
    @{ Kind = "TEXT"; Text = "# 10. Call Days by Designation + DR Coverage % by Full Name" },
    @{ Kind = "TEXT"; Text = "st.subheader(`"10A. Call Days by Designation`")" },
    @{ Kind = "TEXT"; Text = "st.plotly_chart(fig_call_by_designation, use_container_width=True)" },
    @{ Kind = "EMPTY"; Text = "" },
    @{ Kind = "TEXT"; Text = "st.subheader(`"10B. DR Coverage % by Full Name`")" },
    @{ Kind = "TEXT"; Text = "st.plotly_chart(fig_fullname_dr_cov, use_container_width=True)" }
)


# Remember where the original content ended, so we can find the new
# paragraphs by index afterwards (Paragraphs collection is 1-based).
$startCount = $d.Paragraphs.Count
$borderOffset = -1
$i = 0

foreach ($item in $items) {
    $sel.EndKey(6)          # wdStory -> move to end of document
    $sel.TypeParagraph()    # start a brand new paragraph

    if ($item.Kind -eq "TEXT" -or $item.Kind -eq "TEXT_PAGEBREAK") {
        $sel.TypeText($item.Text)
    }
    elseif ($item.Kind -eq "BORDER") {
        $borderOffset = $i
    }

    $i = $i + 1
}

# Now that every paragraph mark has been created, go back and apply the
# bottom border to only the one intended paragraph -- doing this while
# still typing would make every later TypeParagraph() inherit the border.
if ($borderOffset -ge 0) {
    $borderParaIndex = $startCount + $borderOffset + 1
    $borderPara = $d.Paragraphs.Item($borderParaIndex)
    $borders = $borderPara.Borders
    $borders.DistanceFromBottom = 1
    $bottom = $borders.Item(-3)        # wdBorderBottom
    $bottom.LineStyle = 14             # -> thinThickThinMediumGap
    $bottom.LineWidth = 9              # -> w:sz="18"
    $bottom.Color = -16777216          # wdColorAutomatic -> w:color="auto"
}

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
